$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# "update genetics, cgm, eeg": the CGM test-case row's id is refreshed to a
# newly generated CA- id.
$ws.Range("A2").Value = "CA-756V081T"
